$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2
$ws.Range("C2").Value = 1.32
$ws.Range("E2").Value = 1.16

# Row 3
$ws.Range("D3").Value = 1.4

# Row 4
$ws.Range("C4").Value = 1.44
$ws.Range("E4").Value = 1.21
$ws.Range("F4").Value = 1.11

# Row 5
$ws.Range("B5").Value = 1.58

# Row 6
$ws.Range("D6").Value = 1.49
$ws.Range("E6").Value = 1.32

# Row 7
$ws.Range("G7").Value = 1.16

$wb.Save()
